$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.244.93"
$ws.Range("E2").Value = "  -1.31%  "

$ws.Range("D3").Value = "2.274.82"
$ws.Range("E3").Value = "  -1.68%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "112.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.71%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "264.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.91%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.621"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.05%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.607"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.31%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "47.79"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.64%  "

$ws.Range("E11").Value = "  -1.50%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.79"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.03%  "

$ws.Range("E13").Value = "  +0.83%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.46"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.68%  "

$ws.Range("E15").Value = "  -1.71%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.853"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.70%  "

$ws.Range("D17").Value = "2.270.73"
$ws.Range("E17").Value = "  -1.81%  "

$ws.Range("D18").Value = "43.165.25"
$ws.Range("E18").Value = "  -1.36%  "

$ws.Range("E19").Value = "  -2.29%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.76"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.20%  "

$ws.Range("E21").Value = "  -2.04%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.50"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.42%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "231.61"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.69"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.13%  "

$ws.Range("E25").Value = "  -1.16%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.28"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.59%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.92"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.58%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.25"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.92%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.33"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.03%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.24"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.55%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "172.19"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.08%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.29"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.93%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0905"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.84%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.75"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.01%  "

$ws.Range("E36").Value = "  +0.28%  "

$ws.Range("E37").Value = "  -1.81%  "

$ws.Range("E38").Value = "  -1.60%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.83"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.13%  "

$ws.Range("E40").Value = "  -6.69%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.61"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +9.22%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "76.80"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.96%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "13.78"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.44%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.235"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.62%  "

$ws.Range("E45").Value = "  +1.43%  "

$ws.Range("E46").Value = "  +0.11%  "

$ws.Range("E47").Value = "  -1.72%  "

$ws.Range("E48").Value = "  -2.46%  "

$ws.Range("B49").Value = "TrustWalletToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.25"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.39%  "

$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "101.50"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.98%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0992"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.69%  "
